$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds text-formatted numbers (e.g. "319.27", "1.091").
# Force text NumberFormat before assignment so Excel does not silently
# coerce these into numeric values (which would lose formatting like
# trailing zeros, e.g. "0.00001080").
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.002.59"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.868.77"
$ws.Range("E3").Value = "  -2.71%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.27"
$ws.Range("E5").Value = "  -2.40%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5098"
$ws.Range("E7").Value = "  -1.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3926"
$ws.Range("E8").Value = "  -2.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08179"
$ws.Range("E9").Value = "  -3.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.30"
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.091"
$ws.Range("E11").Value = "  -2.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.76"
$ws.Range("E12").Value = "  +2.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.862.19"
$ws.Range("E13").Value = "  -3.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.257"
$ws.Range("E14").Value = "  -1.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.149"
$ws.Range("E15").Value = "  -2.79%  "
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.78"
$ws.Range("E17").Value = "  -4.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001080"
$ws.Range("E18").Value = "  -3.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06383"
$ws.Range("E19").Value = "  -5.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.89"
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.979.48"
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.800"
$ws.Range("E23").Value = "  -4.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.08"
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.156"
$ws.Range("E25").Value = "  -2.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.085.59"
$ws.Range("E26").Value = "  -2.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.12"
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.92"
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.227"
$ws.Range("E29").Value = "  -9.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.93"
$ws.Range("E30").Value = "  -1.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.048"
$ws.Range("E31").Value = "  -2.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1034"
$ws.Range("E32").Value = "  -2.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.901"
$ws.Range("E33").Value = "  -2.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.731"
$ws.Range("E34").Value = "  +1.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02414"
$ws.Range("E35").Value = "  -4.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.245"
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06330"
$ws.Range("E37").Value = "  -3.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2142"
$ws.Range("E38").Value = "  -3.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.172"
$ws.Range("E39").Value = "  -5.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.499"
$ws.Range("E40").Value = "  -5.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6297"
$ws.Range("E41").Value = "  -3.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.205"
$ws.Range("E42").Value = "  -3.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.24"
$ws.Range("E43").Value = "  -2.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.92"
$ws.Range("E45").Value = "  -2.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5883"
$ws.Range("E46").Value = "  -4.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.626"
$ws.Range("E47").Value = "  -3.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.985"
$ws.Range("E48").Value = "  -3.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.58"
$ws.Range("E49").Value = "  -2.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.204"
$ws.Range("E50").Value = "  -3.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.138"
$ws.Range("E51").Value = "  -1.59%  "
